$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A13").Value = "next empty cell"
$ws.Range("A14").Value = "next empty cell"
